$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.211.79"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "2.304.74"
$ws.Range("E3").Value = "  -2.14%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.65"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.01"
$ws.Range("E6").Value = "  -3.30%  "

$ws.Range("E7").Value = "  -1.09%  "

$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.89"
$ws.Range("E10").Value = "  -4.36%  "

$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.41"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("E15").Value = "  -3.21%  "

$ws.Range("D16").Value = "2.650.98"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").Value = "2.307.48"
$ws.Range("E17").Value = "  -5.11%  "

$ws.Range("D18").Value = "42.017.72"
$ws.Range("E18").Value = "  -1.77%  "

$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "289.15"
$ws.Range("E21").Value = "  +12.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.67"
$ws.Range("E22").Value = "  -3.90%  "

$ws.Range("E23").Value = "  -1.67%  "

$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("E25").Value = "  +7.22%  "

$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("E27").Value = "  -4.05%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  +4.48%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.47"
$ws.Range("E29").Value = "  +2.50%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "164.55"
$ws.Range("E30").Value = "  -5.93%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.47"
$ws.Range("E31").Value = "  -3.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0885"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.91"
$ws.Range("E34").Value = "  -3.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("E36").Value = "  -8.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.96"
$ws.Range("E38").Value = "  +10.88%  "

$ws.Range("E39").Value = "  -2.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -5.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.39"
$ws.Range("E41").Value = "  +22.98%  "

$ws.Range("E42").Value = "  +1.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.07"
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("E44").Value = "  -4.60%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.35"
$ws.Range("E46").Value = "  +3.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.14"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.10"
$ws.Range("E48").Value = "  +5.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.35"
$ws.Range("E50").Value = "  -2.76%  "

$ws.Range("E51").Value = "  +1.64%  "
